$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddShift")

# Update "Execute" column values from N to Y for specific rows
$ws.Range("A5").Value = "Y"
$ws.Range("A6").Value = "Y"
$ws.Range("A10").Value = "Y"
$ws.Range("A18").Value = "Y"
$ws.Range("A27").Value = "Y"

# Update ActualResult (J26) with error message
$ws.Range("J26").Value = "เกิดข้อผิดพลาด กรุณาลองใหม่อีกครั้ง !!!"

# Update ActualResult (J27) to match expected result, and Result to Pass
$ws.Range("J27").Value = "กรุณากรอกการรับหน้าที่"
$ws.Range("K27").Value = "Pass"

# Move selection/view to B7 and scroll to top (remove topLeftCell override)
$ws.Range("B7").Select()

$wb.Save()
